$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Agosto sheet: fix one value, append five new expense rows
# ---------------------------------------------------------------------------
$agosto = $wb.Worksheets.Item("Agosto")

# C4 was 120, corrected to 125.51
$agosto.Range("C4").Value = 125.51

# New rows 11-15 (order matches the new shared-string entries)
$agosto.Range("A11").Value = "Pelicula Motorola One"
$agosto.Range("D11").Value = "S91Loja"
$agosto.Range("B11").Value = 39.97

$agosto.Range("A13").Value = "Conta Luz"
$agosto.Range("C13").Value = 250

$agosto.Range("A12").Value = "Conta Telefone"
$agosto.Range("B12").Value = 140

$agosto.Range("A14").Value = "Anel"
$agosto.Range("C14").Value = 15.34
$agosto.Range("D14").Value = "Aliexpress"

$agosto.Range("A15").Value = "Pedidos m" + [char]0x00E3 + "e"
$agosto.Range("B15").Value = 263.8
$agosto.Range("D15").Value = "Aliexpress"

# ---------------------------------------------------------------------------
# Julho sheet: the "Bateria Carregador Wii" amount actually belongs to the
# "Nao gasto" column, not "Valor" - move it from B6 to C6
# ---------------------------------------------------------------------------
$julho = $wb.Worksheets.Item("Julho")

$amount = $julho.Range("B6").Value2
$julho.Range("B6").Clear()
$julho.Range("C6").Value = $amount

$null = $julho.Range("C6").Select()

# Agosto stays the tab-selected sheet - re-activate it and leave its own
# selection on the new last-used cell
$null = $agosto.Activate()
$null = $agosto.Range("B16").Select()

$excel.Calculate()
